$wb = $excel.ActiveWorkbook

# --- Sheet: ALC (index 1) ---
$ws = $wb.Worksheets.Item(1)
$ws.Range("J62").Value = 2067.1667
$ws.Range("H62").Value = 1875.375
$ws.Range("L62").Value = 2067.1667
$ws.Range("N62").Value = -3315.1667
$ws.Range("J65").Value = 2067.1667
$ws.Range("L65").Value = 10335.8335
$ws.Range("H65").Value = 1875.375
$ws.Range("N65").Value = -16575.8335
$ws.Range("L86").Value = 3902
$ws.Range("J86").Value = 3902
$ws.Range("N86").Value = -6148
$ws.Range("H86").Value = 4258.5713
$ws.Range("L89").Value = 19510
$ws.Range("H89").Value = 4258.5713
$ws.Range("J89").Value = 3902
$ws.Range("N89").Value = -30742
$ws.Range("M132").Value = -1044741.41
$ws.Range("H132").Value = 337487.47
$ws.Range("K132").Value = 1047271.41
$ws.Range("I132").Value = 349090.47
$ws.Range("J138").Value = 4197.619
$ws.Range("K138").Value = 9674.700000000001
$ws.Range("M138").Value = -4534.700000000001
$ws.Range("L138").Value = 12592.857
$ws.Range("N138").Value = -22872.857
$ws.Range("I138").Value = 3224.9
$ws.Range("H138").Value = 3883.8386
# --- Sheet: ARM (index 2) ---
$ws = $wb.Worksheets.Item(2)
$ws.Range("I2").Value = 2220.111
$ws.Range("H2").Value = 2452837
$ws.Range("K2").Value = 2220.111
$ws.Range("M2").Value = -2107.111
$ws.Range("M32").Value = -3543.6667
$ws.Range("I32").Value = 3830.6667
$ws.Range("H32").Value = 8725.273999999999
$ws.Range("K32").Value = 3830.6667
$ws.Range("J63").Value = 6766.6665
$ws.Range("K63").Value = 3500
$ws.Range("L63").Value = 6766.6665
$ws.Range("N63").Value = -8138.6665
$ws.Range("I63").Value = 3500
$ws.Range("H63").Value = 5133.3335
$ws.Range("M63").Value = -2814
$ws.Range("K66").Value = 17500
$ws.Range("L66").Value = 33833.3325
$ws.Range("J66").Value = 6766.6665
$ws.Range("I66").Value = 3500
$ws.Range("M66").Value = -14068
$ws.Range("N66").Value = -40697.3325
$ws.Range("H66").Value = 5133.3335
$ws.Range("H97").Value = 590
$ws.Range("I97").Value = 590
$ws.Range("K97").Value = 590
$ws.Range("M97").Value = -94
$ws.Range("N102").Value = $null
$ws.Range("I102").Value = 2200
$ws.Range("J102").Value = 0
$ws.Range("H102").Value = 2200
$ws.Range("L102").Value = 0
$ws.Range("K102").Value = 2200
$ws.Range("M102").Value = -578
$ws.Range("K110").Value = 1044.0769
$ws.Range("I110").Value = 1044.0769
$ws.Range("H110").Value = 1219.6316
$ws.Range("M110").Value = 1000.9231
$ws.Range("H116").Value = 2452837
$ws.Range("K116").Value = 2220.111
$ws.Range("M116").Value = 73.88900000000012
$ws.Range("I116").Value = 2220.111
# --- Sheet: BSM (index 3) ---
$ws = $wb.Worksheets.Item(3)
$ws.Range("I3").Value = 2220.111
$ws.Range("M3").Value = -2106.111
$ws.Range("K3").Value = 2220.111
$ws.Range("H3").Value = 2452837
$ws.Range("L94").Value = 1600
$ws.Range("N94").Value = -2502
$ws.Range("I94").Value = 1392.3334
$ws.Range("K94").Value = 1392.3334
$ws.Range("J94").Value = 1600
$ws.Range("M94").Value = -941.3334
$ws.Range("H94").Value = 1422
$ws.Range("J107").Value = 3715.75
$ws.Range("H107").Value = 3131.3333
$ws.Range("I107").Value = 2463.4285
$ws.Range("M107").Value = -543.4285
$ws.Range("L107").Value = 3715.75
$ws.Range("N107").Value = -7555.75
$ws.Range("K107").Value = 2463.4285
$ws.Range("H135").Value = 58042.855
$ws.Range("N135").Value = -68182.85500000001
$ws.Range("L135").Value = 58042.855
$ws.Range("J135").Value = 58042.855
# --- Sheet: CRP (index 4) ---
$ws = $wb.Worksheets.Item(4)
$ws.Range("M12").Value = -1555
$ws.Range("I12").Value = 1725
$ws.Range("K12").Value = 1725
$ws.Range("H12").Value = 2150
$ws.Range("L16").Value = 1299.8572
$ws.Range("K16").Value = 1051.5714
$ws.Range("I16").Value = 1051.5714
$ws.Range("H16").Value = 1175.7142
$ws.Range("N16").Value = -1873.8572
$ws.Range("J16").Value = 1299.8572
$ws.Range("M16").Value = -764.5714
$ws.Range("J113").Value = 1299.8572
$ws.Range("M113").Value = 1118.4286
$ws.Range("I113").Value = 1051.5714
$ws.Range("K113").Value = 1051.5714
$ws.Range("L113").Value = 1299.8572
$ws.Range("N113").Value = -5639.8572
$ws.Range("H113").Value = 1175.7142
$ws.Range("N125").Value = -59900
$ws.Range("L125").Value = 54980
$ws.Range("J125").Value = 54980
$ws.Range("H125").Value = 54980
$ws.Range("M132").Value = -5561
$ws.Range("H132").Value = 3058.8076
$ws.Range("N132").Value = -22557.9995
$ws.Range("L132").Value = 17497.9995
$ws.Range("K132").Value = 8091
$ws.Range("J132").Value = 5832.6665
$ws.Range("I132").Value = 2697
# --- Sheet: CUL (index 5) ---
$ws = $wb.Worksheets.Item(5)
$ws.Range("M68").Value = -1041.8333
$ws.Range("H68").Value = 899.83527
$ws.Range("J68").Value = 1107.1837
$ws.Range("L68").Value = 3321.5511
$ws.Range("N68").Value = -4943.551100000001
$ws.Range("I68").Value = 617.6111
$ws.Range("K68").Value = 1852.8333
$ws.Range("L71").Value = 9964.6533
$ws.Range("H71").Value = 899.83527
$ws.Range("J71").Value = 1107.1837
$ws.Range("N71").Value = -18076.6533
$ws.Range("I71").Value = 617.6111
$ws.Range("M71").Value = -1502.4999
$ws.Range("K71").Value = 5558.4999
$ws.Range("J92").Value = 0
$ws.Range("I92").Value = 600
$ws.Range("K92").Value = 1800
$ws.Range("H92").Value = 600
$ws.Range("L92").Value = 0
$ws.Range("N92").Value = $null
$ws.Range("M92").Value = -552
$ws.Range("J97").Value = 4068
$ws.Range("H97").Value = 2132.1667
$ws.Range("N97").Value = -13196
$ws.Range("L97").Value = 12204
$ws.Range("I97").Value = 196.33333
$ws.Range("K97").Value = 588.99999
$ws.Range("M97").Value = -92.99999000000003
$ws.Range("J107").Value = 25687.85
$ws.Range("H107").Value = 23560.736
$ws.Range("I107").Value = 21750.426
$ws.Range("M107").Value = -63331.278
$ws.Range("L107").Value = 77063.54999999999
$ws.Range("N107").Value = -80903.54999999999
$ws.Range("K107").Value = 65251.278
$ws.Range("L131").Value = 3129414.18
$ws.Range("N131").Value = -3139494.18
$ws.Range("H131").Value = 1001523.9
$ws.Range("J131").Value = 1043138.06
# --- Sheet: GSM (index 6) ---
$ws = $wb.Worksheets.Item(6)
$ws.Range("H122").Value = 8023.385
$ws.Range("I122").Value = 9130.4
$ws.Range("M122").Value = -24941.2
$ws.Range("K122").Value = 27391.2
# --- Sheet: LTW (index 7) ---
$ws = $wb.Worksheets.Item(7)
$ws.Range("M7").Value = -1602.8572
$ws.Range("J7").Value = 1000
$ws.Range("I7").Value = 1714.8572
$ws.Range("L7").Value = 1000
$ws.Range("N7").Value = -1224
$ws.Range("K7").Value = 1714.8572
$ws.Range("H7").Value = 1625.5
$ws.Range("K40").Value = 2048.5
$ws.Range("J40").Value = 0
$ws.Range("M40").Value = -1912.5
$ws.Range("L40").Value = 0
$ws.Range("I40").Value = 2048.5
$ws.Range("H40").Value = 2048.5
$ws.Range("N40").Value = $null
$ws.Range("M126").Value = -2674.571599999999
$ws.Range("H126").Value = 1625.5
$ws.Range("L126").Value = 3000
$ws.Range("K126").Value = 5144.571599999999
$ws.Range("J126").Value = 1000
$ws.Range("N126").Value = -7940
$ws.Range("I126").Value = 1714.8572
$ws.Range("M132").Value = -14427.0314
$ws.Range("H132").Value = 5259.154
$ws.Range("K132").Value = 16957.0314
$ws.Range("I132").Value = 5652.3438
